# "fixed clearing values for kard" - the three "MetKard[]" (cardinality)
# columns on row 2 (I: testDecimalFieldMetKard[], L: testIntegerFieldMetKard[],
# Q: testStringFieldMetKard[]) were left holding stale pipe-joined string
# placeholders ("88888888.0|88888888.0", "88888888|88888888 ") instead of
# being cleared/reset to the same plain numeric value used by every other
# cell in the row. Re-enter them as numbers so they match the rest of the
# row (and the now-unused shared strings fall out of the workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 88888888
$ws.Range("L2").Value = 88888888
$ws.Range("Q2").Value = 88888888

# Leave the selection on the last touched cell of the row.
$ws.Range("R2").Select()
